# Daily attendance processing - 2025-11-04 20:48:09
# For every row in the "Recorded By" column (G), when the recorded-by list
# contains more than one comma-separated entry and one of those entries is
# the "System"/"system" auto-recorder, move it to the back of the list by
# reversing the order of the comma-separated entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Text

    if ($val -eq $null -or $val -eq "") {
        continue
    }

    $parts = $val -split ", "
    $n = $parts.Length

    if ($n -le 1) {
        continue
    }

    $hasSystem = $val -match "System"
    if (-not $hasSystem) {
        continue
    }

    $result = $parts[$n - 1]
    for ($i = $n - 2; $i -ge 0; $i--) {
        $result = $result + ", " + $parts[$i]
    }

    $cell.Value = $result
}
